$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font + border + center/top alignment) from an
# existing header cell into the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2-46).
$ws.Range("AD2:AD46").Value = 80
$ws.Range("AE2:AE46").Value = 81
$ws.Range("AF2:AF46").Value = 0
